# Mob.xlsx update
# - The "DROP LIST" entries for the Slime-type mobs (rows 2-5) are shortened
#   from "2200|2101|-1" to "2101|-1".
# - The "DROP LIST" entries for the next Slime-type mobs (rows 9-11) are
#   shortened from "2200|2100|-1" to "2100|-1".
# - The active cell selection on the sheet moves to D19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "2101|-1"
$ws.Range("D3").Value2 = "2101|-1"
$ws.Range("D4").Value2 = "2101|-1"
$ws.Range("D5").Value2 = "2101|-1"

$ws.Range("D9").Value2  = "2100|-1"
$ws.Range("D10").Value2 = "2100|-1"
$ws.Range("D11").Value2 = "2100|-1"

# Update the saved selection/active cell to match the author's final position.
[void]$ws.Range("D19").Select()
